# Atualização agentes responsáveis dispensa eletronica
# Adds four new rows (4-7) to Planilha1 with new agentes / postos / funções,
# reusing the same look (font/border/alignment) as the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 3) down onto the
# four new rows so the new cells pick up the same style (border, font,
# alignment) as the rest of the table.
$ws.Range("A3:C3").Copy()
$ws.Range("A4:C7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 4: Thiago Martins Amorim now also listed as "Agente Fiscal"
$ws.Range("B4").Value = "Agente Fiscal"
$ws.Range("C4").Value = "Capitão de Fragata (IM)"
$ws.Range("A4").Value = "THIAGO MARTINS AMORIM"

# Row 5: Guilherme Kirschner de Siqueira Campos - Agente Fiscal Substituto
$ws.Range("B5").Value = "Agente Fiscal Substituto"
$ws.Range("C5").Value = "Capitão de Corveta (IM)"
$ws.Range("A5").Value = "GUILHERME KIRSCHNER DE SIQUEIRA CAMPOS"

# Row 6: header-like row for the credit manager table
$ws.Range("B6").Value = "Gerente de Crédito"
$ws.Range("C6").Value = "Posto-Graduacao"
$ws.Range("A6").Value = "NOME COMPLETO"

# Row 7: Ramon de Lima Fernandes - Encarregado da Divisão de Subsistência
$ws.Range("B7").Value = "Encarregado da Divisão de Subsistência"
$ws.Range("C7").Value = "Capitão de Corveta (IM)"
$ws.Range("A7").Value = "RAMON DE LIMA FERNANDES"

# Match the selection left behind in the authored workbook.
[void]$ws.Range("B11").Select()
